# Update Query of intern
# Row 2 corresponds to intern "Rhythm jayee" (jayeerythm8@gmail.com).
# Update her Role, DescriptionOfWork, StartDate and EndDate values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = " Web developer"
$ws.Range("E2").Value = "Frontend & "
$ws.Range("F2").Value = "2020-06-02T18:30:00.000Z"
$ws.Range("G2").Value = "2020-06-15T18:30:00.000Z"
